{"js": "// Update the \"Date:\" field from 16/04/2014 to 07/08/2014, and update the\n// board \"Size\" from 136.4 x 84.2 mm to 138.4 x 86.2 mm.\n\nconst body = context.document.body;\n\n// 1) Date: 16/04/2014 -> 07/08/2014\nconst dateResults = body.search(\"16/04/2014\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"07/08/2014\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Size: \"136.4 x 84.2 mm\" -> \"138.4 x 86.2 mm\"\n//    The width \"36.4\" becomes \"38.4\" ...\nconst widthResults = body.search(\"36.4\", { matchCase: true });\nwidthResults.load(\"items\");\nawait context.sync();\n\nif (widthResults.items.length > 0) {\n  widthResults.items[0].insertText(\"38.4\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n//    ... and the height \"4.2\" (in \" x 84.2 mm\") becomes \"6.2\".\nconst heightResults = body.search(\"4.2\", { matchCase: true });\nheightResults.load(\"items\");\nawait context.sync();\n\nif (heightResults.items.length > 0) {\n  heightResults.items[0].insertText(\"6.2\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the \"Date:\" field from 16/04/2014 to 07/08/2014, and update the\n# board \"Size\" from 136.4 x 84.2 mm to 138.4 x 86.2 mm.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceAll = 2\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# 1) Date: 16/04/2014 -> 07/08/2014\nReplace-DocText \"16/04/2014\" \"07/08/2014\"\n\n# 2) Size: \"136.4 x 84.2 mm\" -> \"138.4 x 86.2 mm\"\nReplace-DocText \"36.4\" \"38.4\"\nReplace-DocText \"4.2\" \"6.2\"\n"}
